$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Local node (Alerta) text fix: "no Sistema" -> "do Sistema"
$ws.Range("D31").Value = "Exibir Alertas do Sistema."

# Status updates for the "Sensor" / "Alerta" tasks that are now complete
$ws.Range("E30").Value = "Pronto"
$ws.Range("E31").Value = "Pronto"
$ws.Range("E32").Value = "Pronto"

# Restore the view state (scroll position / active selection) as left by the author
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("I26").Select()
